# Update PureThermal pinout diagram: the "VIN (3.6V - 5V)" rounded-rectangle
# label (inside the deeply-nested pinout group on slide 1) is relabeled to
# "Not Connected".
#
# The shape lives several group levels deep (Group 7 > Group 1 > Group 100 >
# Rounded Rectangle 50), but this COM host exposes a flattened GroupItems
# collection off the outermost group, so we just scan it for the matching
# text instead of hard-coding fragile indices.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$oldText = "VIN (3.6V - 5V)"
$newText = "Not Connected"

$replaced = $false

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $top = $s.Shapes.Item($i)

    if ($top.HasTextFrame) {
        if ($top.TextFrame.TextRange.Text -eq $oldText) {
            $top.TextFrame.TextRange.Text = $newText
            $replaced = $true
        }
    }

    if ($top.Type -eq 6) {
        # msoGroup - walk the (flattened) group items looking for the label.
        for ($j = 1; $j -le $top.GroupItems.Count; $j++) {
            $sh = $top.GroupItems.Item($j)
            if ($sh.HasTextFrame) {
                if ($sh.TextFrame.TextRange.Text -eq $oldText) {
                    $sh.TextFrame.TextRange.Text = $newText
                    $replaced = $true
                }
            }
        }
    }
}

Write-Host "Replaced '$oldText' -> '$newText':" $replaced
